$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that can look numeric (e.g. "596.85" or
# thousand-dot-grouped values like "61.272.61"); force text storage so
# Excel does not silently coerce them into floating-point numbers, then
# drop the resulting cell style back to Normal so no stray style index
# is left on the cell (matching the original un-styled inline strings).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.272.61"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.58%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.923.53"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.35"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.08%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.501"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.19%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.94"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.93%  "
$ws.Range("E10").Value = "  -1.81%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.437"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.38%  "
$ws.Range("E12").Value = "  -0.81%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "33.36"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.93%  "
$ws.Range("E14").Value = "  +0.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.407.79"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "61.217.46"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.50%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.923.03"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.08%  "
$ws.Range("E18").Value = "  -0.77%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "431.28"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.51"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.673"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.30%  "
$ws.Range("E22").Value = "  -0.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "81.74"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.56%  "
$ws.Range("E25").Value = "  -2.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.74"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.72%  "
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("E28").Value = "  -4.61%  "
$ws.Range("E29").Value = "  -0.87%  "
$ws.Range("E30").Value = "  -2.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.64"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.84%  "
$ws.Range("E32").Value = "  +1.15%  "
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0877"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.24%  "
$ws.Range("E35").Value = "  -0.21%  "
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("E37").Value = "  -1.66%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.00"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.60%  "
$ws.Range("E39").Value = "  -0.43%  "
$ws.Range("E40").Value = "  -0.23%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "42.44"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.10%  "
$ws.Range("E42").Value = "  -2.37%  "
$ws.Range("E43").Value = "  -0.58%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.689.28"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "133.58"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.98%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "363.66"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.61%  "
$ws.Range("E47").Value = "  +0.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.53"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.76%  "
$ws.Range("E49").Value = "  -1.24%  "
$ws.Range("E50").Value = "  -1.02%  "
$ws.Range("E51").Value = "  -0.59%  "
